# Fix dopo primo deploy e gestione date creazione/modifica su Progetti, Persone e Cost Rate
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rinumerazione codici progetto (la stringa condivisa viene aggiornata per A2/A3/A4) ---
$ws.Range("A2").Value = "THE_888_02"
$ws.Range("A3").Value = "THE_888_03"
$ws.Range("A4").Value = "THE_888_04"

# --- Nuova colonna O: ActivityOn ---
$ws.Range("O1").Value = "ActivityOn"
$ws.Range("O3").Value = "X"

# Applica alla nuova intestazione O1 lo stesso stile delle altre intestazioni (riga 1)
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("O1").Value = "ActivityOn"

$excel.CutCopyMode = 0

# Larghezza colonna O (allineata alla colonna N)
$ws.Columns.Item(15).ColumnWidth = 18.5

# Seleziona la nuova cella attiva
$ws.Range("O2").Select()
